# Fix heat rate modeling syntax
$wb = $excel.ActiveWorkbook

# --- Costs and Revenues (row 2) ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76326.84612799998
$ws.Range("C2").Value = 2200
$ws.Range("D2").Value = 9307.780929750721
$ws.Range("E2").Value = 2375
$ws.Range("F2").Value = 33385.14722637499

# --- Capacities (row 4 - Owned Batteries) ---
$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 163

# --- PV Dispatch ---
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("M2").Value = 93.59999999999999
$ws.Range("I3").Value = 0
$ws.Range("M3").Value = 73.45311702887469
$ws.Range("O3").Value = 72.8
$ws.Range("K4").Value = 39.98312417100301
$ws.Range("M4").Value = 83.2
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("R4").Value = 0

# --- Battery Input ---
$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("M2").Value = 70.2
$ws.Range("I3").Value = 0
$ws.Range("M3").Value = 50.05311702887467
$ws.Range("O3").Value = 72.8
$ws.Range("K4").Value = 39.98312417100301
$ws.Range("M4").Value = 59.8
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("R4").Value = 0

# --- Battery Output ---
$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("S2").Value = 9.048519999999844

# --- State of Charge (rows 2-4, columns B-Y) ---
$ws = $wb.Worksheets.Item("State of Charge")

$ws.Range("B2").Value = 189.4909090909091
$ws.Range("C2").Value = 169.7939393939394
$ws.Range("D2").Value = 156.6626262626263
$ws.Range("E2").Value = 143.5313131313131
$ws.Range("F2").Value = 130.4
$ws.Range("G2").Value = 143.27
$ws.Range("H2").Value = 171.584
$ws.Range("I2").Value = 192.176
$ws.Range("J2").Value = 215.342
$ws.Range("K2").Value = 261.674
$ws.Range("L2").Value = 323.45
$ws.Range("M2").Value = 392.948
$ws.Range("N2").Value = 470.168
$ws.Range("O2").Value = 531.944
$ws.Range("P2").Value = 585.9979999999999
$ws.Range("Q2").Value = 632.3299999999999
$ws.Range("R2").Value = 650.348
$ws.Range("S2").Value = 641.2080808080809
$ws.Range("T2").Value = 609.6929292929294
$ws.Range("U2").Value = 491.5111111111111
$ws.Range("V2").Value = 393.0262626262627
$ws.Range("W2").Value = 314.2383838383839
$ws.Range("X2").Value = 261.7131313131313
$ws.Range("Y2").Value = 222.3191919191919

$ws.Range("B3").Value = 182.9252525252525
$ws.Range("C3").Value = 163.2282828282828
$ws.Range("D3").Value = 150.0969696969697
$ws.Range("E3").Value = 150.0969696969697
$ws.Range("F3").Value = 150.0969696969697
$ws.Range("G3").Value = 130.4
$ws.Range("H3").Value = 130.4
$ws.Range("I3").Value = 130.4
$ws.Range("J3").Value = 192.176
$ws.Range("K3").Value = 192.176
$ws.Range("L3").Value = 284.84
$ws.Range("M3").Value = 334.3925858585859
$ws.Range("N3").Value = 334.3925858585859
$ws.Range("O3").Value = 406.4645858585859
$ws.Range("P3").Value = 429.6305858585859
$ws.Range("Q3").Value = 455.3705858585859
$ws.Range("R3").Value = 486.2585858585859
$ws.Range("S3").Value = 465.2484848484848
$ws.Range("T3").Value = 333.9353535353536
$ws.Range("U3").Value = 333.9353535353536
$ws.Range("V3").Value = 333.9353535353536
$ws.Range("W3").Value = 255.1474747474747
$ws.Range("X3").Value = 255.1474747474747
$ws.Range("Y3").Value = 215.7535353535353

$ws.Range("B4").Value = 169.7939393939394
$ws.Range("C4").Value = 150.0969696969697
$ws.Range("D4").Value = 150.0969696969697
$ws.Range("E4").Value = 150.0969696969697
$ws.Range("F4").Value = 150.0969696969697
$ws.Range("G4").Value = 130.4
$ws.Range("H4").Value = 130.4
$ws.Range("I4").Value = 130.4
$ws.Range("J4").Value = 140.696
$ws.Range("K4").Value = 180.279292929293
$ws.Range("L4").Value = 252.351292929293
$ws.Range("M4").Value = 311.553292929293
$ws.Range("N4").Value = 311.553292929293
$ws.Range("O4").Value = 311.553292929293
$ws.Range("P4").Value = 352.7372929292929
$ws.Range("Q4").Value = 373.3292929292929
$ws.Range("R4").Value = 373.3292929292929
$ws.Range("S4").Value = 373.3292929292929
$ws.Range("T4").Value = 242.0161616161616
$ws.Range("U4").Value = 242.0161616161616
$ws.Range("V4").Value = 242.0161616161616
$ws.Range("W4").Value = 242.0161616161616
$ws.Range("X4").Value = 242.0161616161616
$ws.Range("Y4").Value = 202.6222222222222

# --- Feed in from Type 2 ---
$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("S2").Value = 1.351480000000159
